$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update status for a few tasks that moved from "no comenzado" to "terminado"
$ws.Range("B21").Value = "terminado"
$ws.Range("B33").Value = "terminado"
$ws.Range("B44").Value = "terminado"
$ws.Range("B54").Value = "terminado"

# Drop the existing AutoFilter so it can be rebuilt over the new range/criteria
$ws.AutoFilterMode = $false

# Re-apply the AutoFilter over the full data range, filtering column B (estado)
# down to just "no comenzado" (this also hides/shows all rows appropriately)
$rng = $ws.Range("A1:C60")
$rng.AutoFilter(2, @("no comenzado"), 7)

# Keep the _FilterDatabase defined name in sync with the new filter range
$wb.Names.Item("Hoja1!_FilterDatabase").RefersTo = "=Hoja1!`$A`$1:`$C`$60"

# Update the view: scroll back to the top and move the active selection
$ws.Activate()
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C40").Select()
